$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two brand-new rows before row 148, shifting the old row 148 (and below) down by 2.
$ws.Rows.Item(148).Resize(2).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# --- Row 149 now needs to hold what used to be row 147 (untouched original data) ---
$ws.Cells.Item(149, 1).Value = 3
$ws.Cells.Item(149, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(149, 3).Value = "Coquimbo"
$ws.Cells.Item(149, 4).Value = 44552
$ws.Cells.Item(149, 4).NumberFormat = $ws.Cells.Item(147, 4).NumberFormat
$ws.Cells.Item(149, 5).Value = 5
$ws.Cells.Item(149, 6).Value = "Fruta"
$ws.Cells.Item(149, 7).Value = 100103
$ws.Cells.Item(149, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(149, 9).Value = 100103002
$ws.Cells.Item(149, 10).Value = "Ciruela"
$ws.Cells.Item(149, 11).Value = "Black Amber"
$ws.Cells.Item(149, 12).Value = "Primera"
$ws.Cells.Item(149, 13).Value = 60
$ws.Cells.Item(149, 14).Value = 12000
$ws.Cells.Item(149, 15).Value = 12000
$ws.Cells.Item(149, 16).Value = 12000
$ws.Cells.Item(149, 17).Value = "`$/bandeja 10 kilos granel"
$ws.Cells.Item(149, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(149, 19).Value = 1200
$ws.Cells.Item(149, 20).Value = 10

# --- Update row 147 (Primera) in place with the new weekly figures ---
$ws.Cells.Item(147, 4).Value = 44595
$ws.Cells.Item(147, 13).Value = 75
$ws.Cells.Item(147, 14).Value = 13000
$ws.Cells.Item(147, 15).Value = 13000
$ws.Cells.Item(147, 16).Value = 13000
$ws.Cells.Item(147, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(147, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(147, 19).Value = 867
$ws.Cells.Item(147, 20).Value = 15

# --- Fill newly inserted row 148 (Segunda) with the new weekly figures ---
$ws.Cells.Item(148, 1).Value = 3
$ws.Cells.Item(148, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(148, 3).Value = "Coquimbo"
$ws.Cells.Item(148, 4).Value = 44595
$ws.Cells.Item(148, 4).NumberFormat = $ws.Cells.Item(147, 4).NumberFormat
$ws.Cells.Item(148, 5).Value = 5
$ws.Cells.Item(148, 6).Value = "Fruta"
$ws.Cells.Item(148, 7).Value = 100103
$ws.Cells.Item(148, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(148, 9).Value = 100103002
$ws.Cells.Item(148, 10).Value = "Ciruela"
$ws.Cells.Item(148, 11).Value = "Black Amber"
$ws.Cells.Item(148, 12).Value = "Segunda"
$ws.Cells.Item(148, 13).Value = 70
$ws.Cells.Item(148, 14).Value = 12000
$ws.Cells.Item(148, 15).Value = 12000
$ws.Cells.Item(148, 16).Value = 12000
$ws.Cells.Item(148, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(148, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(148, 19).Value = 800
$ws.Cells.Item(148, 20).Value = 15

$wb.Save()
